$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9811958074569702
$ws.Range("B1").Value = 3.322385549545288
$ws.Range("C1").Value = 4.037312507629395
$ws.Range("D1").Value = 3.067291498184204
$ws.Range("E1").Value = 1.328108668327332
